$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new data rows right before the current row 622, pushing the
# existing rows 622-667 down to 625-670 (dimension grows from T667 to T670).
$ws.Rows.Item(622).Insert()
$ws.Rows.Item(623).Insert()
$ws.Rows.Item(624).Insert()

# New row 622: Artic Star / Primera
$ws.Range("A622").Value = 10
$ws.Range("B622").Value = "Vega Modelo de Temuco"
$ws.Range("C622").Value = "La Araucanía"
$ws.Range("D622").Value = 44931
$ws.Range("E622").Value = 9
$ws.Range("F622").Value = "Fruta"
$ws.Range("G622").Value = 100103
$ws.Range("H622").Value = "Frutos de hueso (carozo)"
$ws.Range("I622").Value = 100103006
$ws.Range("J622").Value = "Nectarín"
$ws.Range("K622").Value = "Artic Star"
$ws.Range("L622").Value = "Primera"
$ws.Range("M622").Value = 185
$ws.Range("N622").Value = 17000
$ws.Range("O622").Value = 17000
$ws.Range("P622").Value = 17000
$ws.Range("Q622").Value = "$/bandeja 18 kilos granel"
$ws.Range("R622").Value = "Región de O'Higgins"
$ws.Range("S622").Value = 944
$ws.Range("T622").Value = 18

# New row 623: Early John / Primera
$ws.Range("A623").Value = 10
$ws.Range("B623").Value = "Vega Modelo de Temuco"
$ws.Range("C623").Value = "La Araucanía"
$ws.Range("D623").Value = 44931
$ws.Range("E623").Value = 9
$ws.Range("F623").Value = "Fruta"
$ws.Range("G623").Value = 100103
$ws.Range("H623").Value = "Frutos de hueso (carozo)"
$ws.Range("I623").Value = 100103006
$ws.Range("J623").Value = "Nectarín"
$ws.Range("K623").Value = "Early John"
$ws.Range("L623").Value = "Primera"
$ws.Range("M623").Value = 185
$ws.Range("N623").Value = 17000
$ws.Range("O623").Value = 17000
$ws.Range("P623").Value = 17000
$ws.Range("Q623").Value = "$/bandeja 18 kilos granel"
$ws.Range("R623").Value = "Región de O'Higgins"
$ws.Range("S623").Value = 944
$ws.Range("T623").Value = 18

# New row 624: Super Queen / Especial
$ws.Range("A624").Value = 10
$ws.Range("B624").Value = "Vega Modelo de Temuco"
$ws.Range("C624").Value = "La Araucanía"
$ws.Range("D624").Value = 44931
$ws.Range("E624").Value = 9
$ws.Range("F624").Value = "Fruta"
$ws.Range("G624").Value = 100103
$ws.Range("H624").Value = "Frutos de hueso (carozo)"
$ws.Range("I624").Value = 100103006
$ws.Range("J624").Value = "Nectarín"
$ws.Range("K624").Value = "Super Queen"
$ws.Range("L624").Value = "Especial"
$ws.Range("M624").Value = 155
$ws.Range("N624").Value = 28000
$ws.Range("O624").Value = 28000
$ws.Range("P624").Value = 28000
$ws.Range("Q624").Value = "$/caja 20 kilos empedrada"
$ws.Range("R624").Value = "Región de O'Higgins"
$ws.Range("S624").Value = 1400
$ws.Range("T624").Value = 20
